# Auto-generated Word COM-interop script to transform LOQ4104.docx
# per the commit diff (paragraph content rotation + Avaliacao run reshuffle).
$d = $word.ActiveDocument

# --- Whole-paragraph content replacements (single-run paragraphs) ---

$d.Paragraphs.Item(6).Range.Text = 'Introdução aos Processos Químicos Industriais; Química Fina; Petroquímica, Fertilizantes; Papel e Celulose; Óleos e Gorduras; Indústria Cerâmica, Ácidos Sulfônicos e Sulfatos Orgânicos; Aminas, e Álcoois.'

$d.Paragraphs.Item(7).Range.Text = 'Introduction to Industrial Chemistry Processes; Fine Chemicals; Petrochemicals, Fertilizers; Paper And Cellulose; Oil and fat; Ceramic Industry, Sulfonic Acids and Organic Sulfates; Amines, and Alcohols.'

$d.Paragraphs.Item(9).Range.Text = 'Proporcionar aos alunos uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria-prima em produto. Serão estudados os processos das indústrias químicas de base, transformação e de base orgânica.'

$d.Paragraphs.Item(11).Range.Text = '1- Introdução aos Processos Químicos Industriais: 1.1- Definição e Objetivos de um Processo Químico, 1.2- Operações Unitárias e Processos Unitários, 1.3- Tipos de Processos, 1.4- Fluxogramas, 1.5- Setores da Industria Química; 2- Petroquímica: 2.1- Visão Geral, 2.2- Cadeia produtiva, 2.3- Segmento Cloro - Soda, 2.4- Gás de Síntese, Metanol e Amônia, 3- Fertilizantes: 3.1- Visão Geral, 3.2- Principais Compostos Químicos Utilizados (Ácidos Sulfúrico, Nítrico, Fosfórico e seus respectivos derivados), 3.3- Cadeia Produtiva; 4- Papel e Celulose: 4.1- Visão Geral, 4.2- Cadeia Produtiva; 5- Óleos e Gorduras: 5.1- Visão Geral, 5.2- Fontes de Obtenção de Óleos e Gorduras, 5.3- Principais Compostos Químicos e Efeitos nas Propriedades Físico-Químicas e Organolépticas, 5.4- Processos para a Obtenção de Derivados Graxos: 5.4.1 Esterificação, 5.4.2- Hidrogenação, 5.4.3- Oxidação; 6-' + [char]11 + 'Indústria Cerâmica: 6.1 Cimento: 6.1.1- Visão Geral, 6.1.2- Cadeia Produtiva, 6.2- Vidro: 6.2.1- Visão Geral, 6.2.2- Cadeia Produtiva; 7- Ácidos Sulfônicos e Sulfatos Orgânicos, 7.1- Processo de Sulfonação e Sulfatação; 8- Aminas e Álcoois, 8.1- Processos Oxo e Amino.'

$d.Paragraphs.Item(12).Range.Text = 'Provide students with a current view of industrial processes that use chemical conversion as a route for transforming raw materials into products. The processes of chemical-based, transformation and organic-based industries will be studied.'

$d.Paragraphs.Item(14).Range.Text = 'Método:' + [char]11 + 'Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos.' + [char]11 + 'Critério:' + [char]11 + 'A nota (NOTA) será composta por uma destas opções: prova em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula. A estas opções será incorporado, para cada aluno, seu respectivo percentual de frequência no cálculo da nota final (NF), conforme a fórmula explicitada abaixo:' + [char]11 + 'NF = NOTA x % FREQ.' + [char]11 + 'Norma de Recuperação:' + [char]11 + 'Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita aplicação de prova escrita de recuperação valendo 10,00 pontos.'

$d.Paragraphs.Item(19).Range.Text = '1285870 - Marcos Villela Barcza'

# --- Paragraph 17 ("Avaliacao" list): swap the three content runs that
#     follow the bold "Metodo:"/"Criterio:"/"Norma de recuperacao:" labels,
#     leaving the bold labels themselves untouched. ---

$p17 = $d.Paragraphs.Item(17)
$pStart = $p17.Range.Start
$pEnd = $p17.Range.End

$fMetodo = $d.Range($pStart, $pEnd)
$fMetodo.Find.Execute("Método: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$metodoLabelEnd = $fMetodo.End

$fCriterio = $d.Range($pStart, $pEnd)
$fCriterio.Find.Execute("Critério: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$criterioLabelStart = $fCriterio.Start
$criterioLabelEnd = $fCriterio.End

$fNorma = $d.Range($pStart, $pEnd)
$fNorma.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$normaLabelStart = $fNorma.Start
$normaLabelEnd = $fNorma.End

$metodoContent = 'A nota (NOTA) será composta por uma destas opções: prova em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula. A estas opções será incorporado, para cada aluno, seu respectivo percentual de frequência no cálculo da nota final (NF), conforme a fórmula explicitada abaixo:' + [char]11 + 'NF = NOTA x % FREQ.' + [char]11 + ''
$criterioContent = 'Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita aplicação de prova escrita de recuperação valendo 10,00 pontos.' + [char]11 + ''
$normaContent = 'Material elaborado pelo docente. Livros:' + [char]11 + 'Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.' + [char]11 + 'Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.' + [char]11 + 'Manual Econômico da Indústria Química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.' + [char]11 + 'Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio' + [char]11 + 'de Janeiro: Editora Guanabara Koogan, 2008, c1997. Revistas:' + [char]11 + 'Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/category/revista/' + [char]11 + 'Petróleo & Energia, São Paulo, SP, v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/' + [char]11 + 'Revista FACTO, Publicação da Associação Brasileira das Indústrias de Química Fina, Biotecnologia e suas Especialidades, Rio de Janeiro, RJ, v. 1, n. 1; Disponível em: http://www.abifina.org.br/facto/' + [char]11 + 'Revista Óleos & Gorduras, disponível em: https://www.editorastilo.com.br/revista-oleos-e-gorduras/'

# Apply back-to-front so earlier offsets stay valid as text lengths change.
$d.Range($normaLabelEnd, $pEnd).Text = $normaContent
$d.Range($criterioLabelEnd, $normaLabelStart).Text = $criterioContent
$d.Range($metodoLabelEnd, $criterioLabelStart).Text = $metodoContent

Write-Host "Edit complete"
